# Resolved encapsulation issues in file manager. Name mangling
$wb = $excel.ActiveWorkbook

# --- Sheet "vocabluary": add new vocab row 102 -----------------------------
$wsVocab = $wb.Worksheets.Item("vocabluary")
$wsVocab.Range("A102").Value = "socks"
$wsVocab.Range("B102").Value = "skarpety"
$wsVocab.Range("C102").Value = 10

# --- Sheet "categories": add new category row 11 ---------------------------
$wsCat = $wb.Worksheets.Item("categories")
$wsCat.Range("A11").Value = 10
$wsCat.Range("B11").Value = "Test level a1"

# --- Sheet "test results": fix F11, append rows 12 and 13 ------------------
$wsRes = $wb.Worksheets.Item("test results")

# F11 used to be stored as text "73.46"; it becomes a genuine number.
$wsRes.Range("F11").Value = 73.45999999999999

# New row 12 (test_id 11, "Gienek")
$wsRes.Range("A12").Value = 11
$wsRes.Range("B12").Value = "Gienek"
$wsRes.Range("C12").Value = "18-01-2025 01:32:58"
$wsRes.Range("D12").Value = "EN->PL"
$wsRes.Range("E12").Value = 60
$wsRes.Range("F12").Value = 6.75
$wsRes.Range("G12").Value = "(0/<bound method NewTest.__get_questions_amount of <test.test.NewTest object at 0x104d88980>>)"
$wsRes.Range("H12").Value = "'0.00%"
$wsRes.Range("I12").Value = "Monitory, Keyboard, Test level a1"

# New row 13 (test_id 12, "GIenek") -- F13 stays textual ("7.67"), unlike F12
$wsRes.Range("A13").Value = 12
$wsRes.Range("B13").Value = "GIenek"
$wsRes.Range("C13").Value = "18-01-2025 01:37:01"
$wsRes.Range("D13").Value = "EN->PL"
$wsRes.Range("E13").Value = 60
$wsRes.Range("F13").Value = "'7.67"
$wsRes.Range("G13").Value = "(1/<bound method NewTest.__get_questions_amount of <test.test.NewTest object at 0x105224830>>)"
$wsRes.Range("H13").Value = "'100.00%"
$wsRes.Range("I13").Value = "Monitory, Keyboard, Test level a1"
